$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, whether it is a numeric-looking
# text string in column D that Excel would otherwise silently coerce to
# a Number (dropping significant trailing zeros / using sci. notation).
# Those are written through a Text-formatted cell, then restyled back to
# "Normal" so the saved style index matches the original (unstyled) cell.
$updates = @(
    @{Cell="D2"; Value="26.470.96"; ForceText=$false},
    @{Cell="E2"; Value="  -0.64%  "; ForceText=$false},
    @{Cell="D3"; Value="1.837.18"; ForceText=$false},
    @{Cell="E3"; Value="  -0.89%  "; ForceText=$false},
    @{Cell="D4"; Value="1.000"; ForceText=$true},
    @{Cell="E4"; Value="  -0.05%  "; ForceText=$false},
    @{Cell="D5"; Value="261.25"; ForceText=$true},
    @{Cell="E5"; Value="  -1.35%  "; ForceText=$false},
    @{Cell="E6"; Value="  +0.02%  "; ForceText=$false},
    @{Cell="D7"; Value="0.5375"; ForceText=$true},
    @{Cell="E7"; Value="  +2.03%  "; ForceText=$false},
    @{Cell="D8"; Value="0.3018"; ForceText=$true},
    @{Cell="E8"; Value="  -7.27%  "; ForceText=$false},
    @{Cell="D9"; Value="0.06869"; ForceText=$true},
    @{Cell="E9"; Value="  +1.03%  "; ForceText=$false},
    @{Cell="D10"; Value="17.61"; ForceText=$true},
    @{Cell="E10"; Value="  -7.55%  "; ForceText=$false},
    @{Cell="B11"; Value="Polygon"; ForceText=$false},
    @{Cell="C11"; Value="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; ForceText=$false},
    @{Cell="D11"; Value="0.7369"; ForceText=$true},
    @{Cell="E11"; Value="  -5.90%  "; ForceText=$false},
    @{Cell="B12"; Value="WrappedEther"; ForceText=$false},
    @{Cell="C12"; Value="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; ForceText=$false},
    @{Cell="D12"; Value="1.843.08"; ForceText=$false},
    @{Cell="E12"; Value="  -0.28%  "; ForceText=$false},
    @{Cell="D13"; Value="0.07231"; ForceText=$true},
    @{Cell="E13"; Value="  -7.24%  "; ForceText=$false},
    @{Cell="D14"; Value="89.40"; ForceText=$true},
    @{Cell="E14"; Value="  +0.87%  "; ForceText=$false},
    @{Cell="D15"; Value="4.983"; ForceText=$true},
    @{Cell="E15"; Value="  -0.91%  "; ForceText=$false},
    @{Cell="E16"; Value="  -0.14%  "; ForceText=$false},
    @{Cell="D17"; Value="13.81"; ForceText=$true},
    @{Cell="E17"; Value="  -1.56%  "; ForceText=$false},
    @{Cell="E18"; Value="  +0.01%  "; ForceText=$false},
    @{Cell="D19"; Value="0.000007875"; ForceText=$true},
    @{Cell="E19"; Value="  -1.28%  "; ForceText=$false},
    @{Cell="D20"; Value="26.496.36"; ForceText=$false},
    @{Cell="E20"; Value="  -0.60%  "; ForceText=$false},
    @{Cell="D21"; Value="2.081.40"; ForceText=$false},
    @{Cell="E21"; Value="  +0.25%  "; ForceText=$false},
    @{Cell="D22"; Value="4.587"; ForceText=$true},
    @{Cell="E22"; Value="  -1.20%  "; ForceText=$false},
    @{Cell="D23"; Value="5.963"; ForceText=$true},
    @{Cell="E23"; Value="  -0.78%  "; ForceText=$false},
    @{Cell="E24"; Value="  -3.12%  "; ForceText=$false},
    @{Cell="D25"; Value="142.77"; ForceText=$true},
    @{Cell="E25"; Value="  -0.09%  "; ForceText=$false},
    @{Cell="D26"; Value="2.195"; ForceText=$true},
    @{Cell="E26"; Value="  +0.50%  "; ForceText=$false},
    @{Cell="D27"; Value="1.686"; ForceText=$true},
    @{Cell="E27"; Value="  +0.14%  "; ForceText=$false},
    @{Cell="D28"; Value="16.95"; ForceText=$true},
    @{Cell="E28"; Value="  -0.63%  "; ForceText=$false},
    @{Cell="D29"; Value="110.65"; ForceText=$true},
    @{Cell="E29"; Value="  -1.30%  "; ForceText=$false},
    @{Cell="D30"; Value="4.221"; ForceText=$true},
    @{Cell="E30"; Value="  +0.47%  "; ForceText=$false},
    @{Cell="D31"; Value="0.08823"; ForceText=$true},
    @{Cell="E31"; Value="  +1.03%  "; ForceText=$false},
    @{Cell="D32"; Value="4.018"; ForceText=$true},
    @{Cell="E32"; Value="  -2.45%  "; ForceText=$false},
    @{Cell="D33"; Value="0.04807"; ForceText=$true},
    @{Cell="E33"; Value="  -0.87%  "; ForceText=$false},
    @{Cell="D34"; Value="2.916"; ForceText=$true},
    @{Cell="E34"; Value="  +1.22%  "; ForceText=$false},
    @{Cell="D35"; Value="0.7273"; ForceText=$true},
    @{Cell="E35"; Value="  +0.56%  "; ForceText=$false},
    @{Cell="D36"; Value="1.130"; ForceText=$true},
    @{Cell="E36"; Value="  -0.24%  "; ForceText=$false},
    @{Cell="E37"; Value="  -0.72%  "; ForceText=$false},
    @{Cell="D38"; Value="2.267"; ForceText=$true},
    @{Cell="E38"; Value="  -0.13%  "; ForceText=$false},
    @{Cell="D39"; Value="0.01705"; ForceText=$true},
    @{Cell="D40"; Value="0.4708"; ForceText=$true},
    @{Cell="E40"; Value="  -3.63%  "; ForceText=$false},
    @{Cell="D41"; Value="0.9034"; ForceText=$true},
    @{Cell="E41"; Value="  -0.22%  "; ForceText=$false},
    @{Cell="D42"; Value="107.70"; ForceText=$true},
    @{Cell="E42"; Value="  -2.93%  "; ForceText=$false},
    @{Cell="D43"; Value="5.880"; ForceText=$true},
    @{Cell="E43"; Value="  -1.70%  "; ForceText=$false},
    @{Cell="D44"; Value="1.000"; ForceText=$true},
    @{Cell="E44"; Value="  -0.01%  "; ForceText=$false},
    @{Cell="D45"; Value="7.359"; ForceText=$true},
    @{Cell="E45"; Value="  -4.37%  "; ForceText=$false},
    @{Cell="D46"; Value="9.003"; ForceText=$true},
    @{Cell="E46"; Value="  -0.70%  "; ForceText=$false},
    @{Cell="D47"; Value="0.1239"; ForceText=$true},
    @{Cell="E47"; Value="  +0.07%  "; ForceText=$false},
    @{Cell="D48"; Value="0.4051"; ForceText=$true},
    @{Cell="E48"; Value="  -3.78%  "; ForceText=$false},
    @{Cell="D49"; Value="34.73"; ForceText=$true},
    @{Cell="E49"; Value="  -1.12%  "; ForceText=$false},
    @{Cell="D50"; Value="0.05772"; ForceText=$true},
    @{Cell="E50"; Value="  -2.03%  "; ForceText=$false},
    @{Cell="D51"; Value="0.8901"; ForceText=$true},
    @{Cell="E51"; Value="  +0.00%  "; ForceText=$false}
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $c.NumberFormat = "@"
        $c.Value = $u.Value
        $c.Style = "Normal"
    } else {
        $c.Value = $u.Value
    }
}
